$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.289385795593262
$ws.Range("B1").Value = 1.742064833641052
$ws.Range("C1").Value = 1.660213828086853
$ws.Range("D1").Value = 4.554900169372559
$ws.Range("E1").Value = 1.351874947547913
